$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update "Version" property value (row 3, column B): 1.0.0 -> 2.0.2
$ws.Cells.Item(3, 2).Value = "2.0.2"

# Update "Date" property value (row 8, column B)
$ws.Cells.Item(8, 2).Value = "2025-02-13T16:11:24+00:00"

# Insert a new row after "Contact" (row 10) for the new "Jurisdiction" property
$ws.Rows.Item(11).Insert()

# The freshly inserted row doesn't inherit the table's border/alignment formatting;
# copy that formatting down from the row now below it (old row 11, "Description").
$ws.Range("A12:B12").Copy()
$ws.Range("A11:B11").PasteSpecial(-4122)

$ws.Cells.Item(11, 1).Value = "Jurisdiction"
$ws.Cells.Item(11, 2).Value = ""
